$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new date columns (Jun_27, Jun_26, Jun_26) before the existing
# B:D columns, pushing the old Jun_17/Jun_15/Jun_13/Jun_10 columns right.
$ws.Columns("B:D").Insert()

# New header cells for the inserted date columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the rating placeholder ("UN") for all existing watchlist rows in the
# newly inserted columns (B:D) as well as the columns that the old B:D data
# shifted into (E:G), matching the rest of the table.
$ws.Range("B2:G27").Value = "UN"

# Add the two new watchlist entries at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
